$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "Burger Peach"
$ws.Cells.Item(2, 3).Value = 43
$ws.Cells.Item(2, 4).Value = 1
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 25029
$ws.Cells.Item(2, 7).Value = 1154
$ws.Cells.Item(2, 8).Value = 65.06999999999999
$ws.Cells.Item(2, 9).Value = 46
$ws.Cells.Item(2, 10).Value = 2

$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "Nick Fitzpatrick"
$ws.Cells.Item(3, 3).Value = 30
$ws.Cells.Item(3, 4).Value = 4
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 16869
$ws.Cells.Item(3, 7).Value = 799
$ws.Cells.Item(3, 8).Value = 63.34
$ws.Cells.Item(3, 9).Value = 34
$ws.Cells.Item(3, 10).Value = 1

$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "Yannick den Daggelder"
$ws.Cells.Item(4, 3).Value = 20
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 11319
$ws.Cells.Item(4, 7).Value = 619
$ws.Cells.Item(4, 8).Value = 54.86
$ws.Cells.Item(4, 9).Value = 20
$ws.Cells.Item(4, 10).Value = 1

$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "Niels van Dommelen"
$ws.Cells.Item(5, 3).Value = 19
$ws.Cells.Item(5, 4).Value = 0
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 15894
$ws.Cells.Item(5, 7).Value = 910
$ws.Cells.Item(5, 8).Value = 52.4
$ws.Cells.Item(5, 9).Value = 19
$ws.Cells.Item(5, 10).Value = 0

$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "Rocky Van Den Eeckhoudt"
$ws.Cells.Item(6, 3).Value = 13
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 11777
$ws.Cells.Item(6, 7).Value = 619
$ws.Cells.Item(6, 8).Value = 57.08
$ws.Cells.Item(6, 9).Value = 14
$ws.Cells.Item(6, 10).Value = 0

$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "Lukas G"
$ws.Cells.Item(7, 3).Value = 11
$ws.Cells.Item(7, 4).Value = 1
$ws.Cells.Item(7, 5).Value = 0
$ws.Cells.Item(7, 6).Value = 9785
$ws.Cells.Item(7, 7).Value = 545
$ws.Cells.Item(7, 8).Value = 53.86
$ws.Cells.Item(7, 9).Value = 12
$ws.Cells.Item(7, 10).Value = 0

$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "Nigel Riedel"
$ws.Cells.Item(8, 3).Value = 10
$ws.Cells.Item(8, 4).Value = 1
$ws.Cells.Item(8, 5).Value = 0
$ws.Cells.Item(8, 6).Value = 9077
$ws.Cells.Item(8, 7).Value = 541
$ws.Cells.Item(8, 8).Value = 50.33
$ws.Cells.Item(8, 9).Value = 11
$ws.Cells.Item(8, 10).Value = 0

$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "Sion Foulkes"
$ws.Cells.Item(9, 3).Value = 8
$ws.Cells.Item(9, 4).Value = 0
$ws.Cells.Item(9, 5).Value = 0
$ws.Cells.Item(9, 6).Value = 6951
$ws.Cells.Item(9, 7).Value = 454
$ws.Cells.Item(9, 8).Value = 45.93
$ws.Cells.Item(9, 9).Value = 8
$ws.Cells.Item(9, 10).Value = 0

$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "Noah B"
$ws.Cells.Item(10, 3).Value = 6
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(10, 5).Value = 0
$ws.Cells.Item(10, 6).Value = 6147
$ws.Cells.Item(10, 7).Value = 400
$ws.Cells.Item(10, 8).Value = 46.1
$ws.Cells.Item(10, 9).Value = 6
$ws.Cells.Item(10, 10).Value = 0

$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "Aelbrecht Wesley"
$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(11, 4).Value = 0
$ws.Cells.Item(11, 5).Value = 0
$ws.Cells.Item(11, 6).Value = 2740
$ws.Cells.Item(11, 7).Value = 200
$ws.Cells.Item(11, 8).Value = 41.1
$ws.Cells.Item(11, 9).Value = 5
$ws.Cells.Item(11, 10).Value = 0

$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "Constantinos Mavroudis"
$ws.Cells.Item(12, 3).Value = 4
$ws.Cells.Item(12, 4).Value = 0
$ws.Cells.Item(12, 5).Value = 0
$ws.Cells.Item(12, 6).Value = 1991
$ws.Cells.Item(12, 7).Value = 134
$ws.Cells.Item(12, 8).Value = 44.57
$ws.Cells.Item(12, 9).Value = 4
$ws.Cells.Item(12, 10).Value = 0

$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = "Danny Littler"
$ws.Cells.Item(13, 3).Value = 3
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = ""
$ws.Cells.Item(13, 9).Value = 3
$ws.Cells.Item(13, 10).Value = 0

$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "Tristan Snoep"
$ws.Cells.Item(14, 3).Value = 3
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(14, 6).Value = 1193
$ws.Cells.Item(14, 7).Value = 96
$ws.Cells.Item(14, 8).Value = 37.28
$ws.Cells.Item(14, 9).Value = 3
$ws.Cells.Item(14, 10).Value = 0

$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = "Arnold Van Der Vlies"
$ws.Cells.Item(15, 3).Value = 2
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(15, 6).Value = 2362
$ws.Cells.Item(15, 7).Value = 167
$ws.Cells.Item(15, 8).Value = 42.43
$ws.Cells.Item(15, 9).Value = 2
$ws.Cells.Item(15, 10).Value = 0

$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "Quintin Marais"
$ws.Cells.Item(16, 3).Value = 2
$ws.Cells.Item(16, 4).Value = 0
$ws.Cells.Item(16, 5).Value = 0
$ws.Cells.Item(16, 6).Value = 2388
$ws.Cells.Item(16, 7).Value = 107
$ws.Cells.Item(16, 8).Value = 66.95
$ws.Cells.Item(16, 9).Value = 2
$ws.Cells.Item(16, 10).Value = 0
